$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DomainDelivery_DomainID")
$ws2 = $wb.Worksheets.Item("DomainConversion_DomainID")

# Update the display-level labels (dropping the redundant "Domain " prefix)
$ws1.Range("B3").Value = "Click Based Conversions"
$ws1.Range("B4").Value = "View Based Conversions"

$ws2.Range("B3").Value = "Click Based Conversions"
$ws2.Range("B4").Value = "View Based Conversions"

# Reflect the author's final cursor position on each sheet
$ws1.Activate() | Out-Null
$ws1.Range("C4").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C4").Select() | Out-Null
